$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) with the same header style as the other
# header cells (e.g. G1: bold, bordered, centered) by copying G1's format.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Add the data value for the new column in row 2.
$ws.Range("H2").Value = 1
